$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.360.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.76%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.794.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.58%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'  -0.01%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'307.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.05%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4516"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.06%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3598"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.99%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'46.44"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.55%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.07082"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.91%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.8849"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.53%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07747"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.04%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'19.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.21%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.781.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.72%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'5.283"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.44%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'6.319"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.81%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'84.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.16%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'1.007"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.000008518"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.01%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -0.05%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'14.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.21%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'26.376.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.79%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'4.970"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.35%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +1.08%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'1.996.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.49%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'1.972"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.02%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'151.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.02%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'17.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.50%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.031"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +4.47%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'111.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.34%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'4.844"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.17%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.08686"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.93%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.083"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.72%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'2.744"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +8.73%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'4.443"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.61%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.7214"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -3.41%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.104"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.14%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +0.03%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -1.31%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.01930"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.06%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.05087"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.44%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.855"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.19%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.5061"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.12%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'6.824"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.24%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.1516"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'8.008"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.20%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.005"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.05%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.4630"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.00%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'EnergySwap"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'9.912"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.85%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'Quant"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'100.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.53%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.564"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.63%  "
$ws.Range("E51").Style = "Normal"
